# Rename the "IV250SD" worksheet to "IV250"
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("IV250SD")
$ws.Name = "IV250"

# Update the sheet's remembered selection: was I1 (I1:I1048576) -> now S22 (S22)
$ws.Activate()
$ws.Range("S22").Select() | Out-Null
